$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "as of" date in A1 (keeps its existing style/number format).
$ws.Range("A1").Value = 43937

# New daily rows appended below the existing data (rows 75-81), matching
# the style (date number format) used by the existing date column cells.
$newRows = @(
    @{Row=75; Date=43931; Count=362},
    @{Row=76; Date=43932; Count=276},
    @{Row=77; Date=43933; Count=250},
    @{Row=78; Date=43934; Count=303},
    @{Row=79; Date=43935; Count=289},
    @{Row=80; Date=43936; Count=406},
    @{Row=81; Date=43937; Count=329}
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    # Copy the formatting (style) of the prior row's date cell down, then
    # overwrite the value - this keeps the same numFmt (m/d;@) style used
    # by all other date cells in column A.
    $ws.Range("A" + ($r - 1)).Copy($ws.Range("A" + $r))
    $ws.Range("A" + $r).Value = $entry.Date
    $ws.Range("B" + $r).Value = $entry.Count
}
